$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) to short column codes
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# Title-case connector words (de/del/el/la/las/los/y -> De/Del/El/La/Las/Los/Y)
# in municipality / state name text cells
$ws.Range("B5").Value = 'Pabellón De Arteaga'
$ws.Range("B6").Value = 'Rincón De Romos'
$ws.Range("B10").Value = 'Playas De Rosarito'
$ws.Range("B27").Value = 'Amatenango De La Frontera'
$ws.Range("B30").Value = 'Bejucal De Ocampo'
$ws.Range("B32").Value = 'Benemérito De Las Américas'
$ws.Range("B37").Value = 'Chiapa De Corzo'
$ws.Range("B41").Value = 'Comitán De Domínguez'
$ws.Range("B58").Value = 'Mazapa De Madero'
$ws.Range("B70").Value = 'Salto De Agua'
$ws.Range("B71").Value = 'San Cristóbal De Las Casas'
$ws.Range("B98").Value = 'Guadalupe Y Calvo'
$ws.Range("B99").Value = 'Hidalgo Del Parral'
$ws.Range("B117").Value = 'San Juan De Sabinas'
$ws.Range("B125").Value = 'Villa De Álvarez'
$ws.Range("A127").Value = 'Ciudad De México'
$ws.Range("B143").Value = 'Coneto De Comonfort'
$ws.Range("B154").Value = 'Nombre De Dios'
$ws.Range("B157").Value = 'San Juan De Guadalupe'
$ws.Range("A163").Value = 'Estado De México'
$ws.Range("B163").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B166").Value = 'Almoloya De Juárez'
$ws.Range("B171").Value = 'Atizapán De Zaragoza'
$ws.Range("B175").Value = 'Chapa De Mota'
$ws.Range("B183").Value = 'Ecatepec De Morelos'
$ws.Range("B188").Value = 'Ixtapan De La Sal'
$ws.Range("B189").Value = 'Ixtapan Del Oro'
$ws.Range("B200").Value = 'Naucalpan De Juárez'
$ws.Range("B208").Value = 'San Felipe Del Progreso'
$ws.Range("B209").Value = 'San Simón De Guerrero'
$ws.Range("B210").Value = 'Soyaniquilpan De Juárez'
$ws.Range("B217").Value = 'Tenango Del Valle'
$ws.Range("B224").Value = 'Tlalnepantla De Baz'
$ws.Range("B228").Value = 'Valle De Bravo'
$ws.Range("B229").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B230").Value = 'Villa De Allende'
$ws.Range("B231").Value = 'Villa Del Carbón'
$ws.Range("B242").Value = 'San Miguel De Allende'
$ws.Range("B243").Value = 'Apaseo El Alto'
$ws.Range("B244").Value = 'Apaseo El Grande'
$ws.Range("B251").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B254").Value = 'Jaral Del Progreso'
$ws.Range("B262").Value = 'Purísima Del Rincón'
$ws.Range("B265").Value = 'San Diego De La Unión'
$ws.Range("B267").Value = 'San Francisco Del Rincón'
$ws.Range("B269").Value = 'San Luis De La Paz'
$ws.Range("B270").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B273").Value = 'Valle De Santiago'
$ws.Range("B279").Value = 'Acapulco De Juárez'
$ws.Range("B281").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B282").Value = 'Alcozauca De Guerrero'
$ws.Range("B285").Value = 'Atenango Del Río'
$ws.Range("B287").Value = 'Atoyac De Álvarez'
$ws.Range("B288").Value = 'Ayutla De Los Libres'
$ws.Range("B291").Value = 'Buenavista De Cuéllar'
$ws.Range("B292").Value = 'Chilapa De Álvarez'
$ws.Range("B293").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B294").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B298").Value = 'Coyuca De Benítez'
$ws.Range("B299").Value = 'Coyuca De Catalán'
$ws.Range("B302").Value = 'Cutzamala De Pinzón'
$ws.Range("B307").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B308").Value = 'Iguala De La Independencia'
$ws.Range("B310").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B311").Value = 'Zihuatanejo De Azueta'
$ws.Range("B313").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B326").Value = 'Taxco De Alarcón'
$ws.Range("B328").Value = 'Técpan De Galeana'
$ws.Range("B330").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B331").Value = 'Tixtla De Guerrero'
$ws.Range("B335").Value = 'Tlapa De Comonfort'
$ws.Range("B347").Value = 'Agua Blanca De Iturbide'
$ws.Range("B353").Value = 'Atotonilco De Tula'
$ws.Range("B354").Value = 'Atotonilco El Grande'
$ws.Range("B359").Value = 'Cuautepec De Hinojosa'
$ws.Range("B363").Value = 'Huasca De Ocampo'
$ws.Range("B365").Value = 'Huejutla De Reyes'
$ws.Range("B368").Value = 'Jacala De Ledezma'
$ws.Range("B374").Value = 'Mineral Del Monte'
$ws.Range("B375").Value = 'Mixquiahuala De Juárez'
$ws.Range("B376").Value = 'Molango De Escamilla'
$ws.Range("B378").Value = 'Nopala De Villagrán'
$ws.Range("B379").Value = 'Omitlán De Juárez'
$ws.Range("B380").Value = 'Pachuca De Soto'
$ws.Range("B383").Value = 'Progreso De Obregón'
$ws.Range("B388").Value = 'Santiago De Anaya'
$ws.Range("B392").Value = 'Tenango De Doria'
$ws.Range("B394").Value = 'Tepehuacán De Guerrero'
$ws.Range("B395").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B396").Value = 'Tezontepec De Aldama'
$ws.Range("B403").Value = 'Tula De Allende'
$ws.Range("B404").Value = 'Tulancingo De Bravo'
$ws.Range("B408").Value = 'Zacualtipán De Ángeles'
$ws.Range("B409").Value = 'Zapotlán De Juárez'
$ws.Range("B414").Value = 'Acatlán De Juárez'
$ws.Range("B415").Value = 'Ahualulco De Mercado'
$ws.Range("B419").Value = 'Atotonilco El Alto'
$ws.Range("B424").Value = 'Encarnación De Díaz'
$ws.Range("B428").Value = 'Huejuquilla El Alto'
$ws.Range("B429").Value = 'Ixtlahuacán Del Río'
$ws.Range("B434").Value = 'Lagos De Moreno'
$ws.Range("B441").Value = 'San Juan De Los Lagos'
$ws.Range("B442").Value = 'San Juanito De Escobedo'
$ws.Range("B443").Value = 'San Miguel El Alto'
$ws.Range("B444").Value = 'Santa María De Los Ángeles'
$ws.Range("B446").Value = 'Tamazula De Gordiano'
$ws.Range("B448").Value = 'Tepatitlán De Morelos'
$ws.Range("B449").Value = 'Tizapán El Alto'
$ws.Range("B450").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B457").Value = 'Unión De Tula'
$ws.Range("B516").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B535").Value = 'Coatlán Del Río'
$ws.Range("B541").Value = 'Jonacatepec De Leandro Valle'
$ws.Range("B545").Value = 'Puente De Ixtla'
$ws.Range("B550").Value = 'Tetela Del Volcán'
$ws.Range("B558").Value = 'Bahía De Banderas'
$ws.Range("B561").Value = 'Ixtlán Del Río'
$ws.Range("B566").Value = 'Santa María Del Oro'
$ws.Range("B583").Value = 'San Nicolás De Los Garza'
$ws.Range("B587").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B590").Value = 'Ayoquezco De Aldama'
$ws.Range("B596").Value = 'Constancia Del Rosario'
$ws.Range("B599").Value = 'Cuyamecalco Villa De Zaragoza'
$ws.Range("B601").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B602").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B603").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B604").Value = 'Huautla De Jiménez'
$ws.Range("B605").Value = 'Ixtlán De Juárez'
$ws.Range("B606").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B609").Value = 'Mariscala De Juárez'
$ws.Range("B612").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B613").Value = 'Mixistlán De La Reforma'
$ws.Range("B615").Value = 'Oaxaca De Juárez'
$ws.Range("B616").Value = 'Ocotlán De Morelos'
$ws.Range("B617").Value = 'Pinotepa De Don Luis'
$ws.Range("B619").Value = 'Putla Villa De Guerrero'
$ws.Range("B632").Value = 'San Dionisio Del Mar'
$ws.Range("B634").Value = 'San Felipe Jalapa De Díaz'
$ws.Range("B645").Value = 'San Juan Bautista Lo De Soto'
$ws.Range("B679").Value = 'San Miguel Del Puerto'
$ws.Range("B690").Value = 'San Pedro El Alto'
$ws.Range("B705").Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range("B714").Value = 'Santa Ana Del Valle'
$ws.Range("B723").Value = 'Santa Lucía Del Camino'
$ws.Range("B773").Value = 'Santo Domingo De Morelos'
$ws.Range("B786").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B787").Value = 'Tataltepec De Valdés'
$ws.Range("B789").Value = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
$ws.Range("B790").Value = 'Tlacolula De Matamoros'
$ws.Range("B792").Value = 'Tlalixtac De Cabrera'
$ws.Range("B794").Value = 'Villa De Etla'
$ws.Range("B795").Value = 'Villa De Tututepec'
$ws.Range("B796").Value = 'Villa De Zaachila'
$ws.Range("B798").Value = 'Villa Sola De Vega'
$ws.Range("B799").Value = 'Villa Talea De Castro'
$ws.Range("B800").Value = 'Zimatlán De Álvarez'
$ws.Range("B813").Value = 'Ayotoxco De Guerrero'
$ws.Range("B827").Value = 'Cuetzalan Del Progreso'
$ws.Range("B834").Value = 'Huehuetlán El Chico'
$ws.Range("B835").Value = 'Huehuetlán El Grande'
$ws.Range("B839").Value = 'Huitzilan De Serdán'
$ws.Range("B841").Value = 'Izúcar De Matamoros'
$ws.Range("B848").Value = 'Los Reyes De Juárez'
$ws.Range("B853").Value = 'Palmar De Bravo'
$ws.Range("B865").Value = 'San Salvador El Seco'
$ws.Range("B876").Value = 'Tetela De Ocampo'
$ws.Range("B877").Value = 'Teteles De Avila Castillo'
$ws.Range("B904").Value = 'Amealco De Bonfil'
$ws.Range("B906").Value = 'Cadereyta De Montes'
$ws.Range("B910").Value = 'Jalpan De Serra'
$ws.Range("B911").Value = 'Landa De Matamoros'
$ws.Range("B913").Value = 'Pinal De Amoles'
$ws.Range("B915").Value = 'San Juan Del Río'
$ws.Range("B924").Value = 'Axtla De Terrazas'
$ws.Range("B930").Value = 'Ciudad Del Maíz'
$ws.Range("B940").Value = 'Mexquitic De Carmona'
$ws.Range("B945").Value = 'San Ciro De Acosta'
$ws.Range("B949").Value = 'Santa María Del Río'
$ws.Range("B950").Value = 'Soledad De Graciano Sánchez'
$ws.Range("B959").Value = 'Villa De Guadalupe'
$ws.Range("B960").Value = 'Villa De Ramos'
$ws.Range("B961").Value = 'Villa De Reyes'
$ws.Range("B996").Value = 'Jalpa De Méndez'
$ws.Range("B1021").Value = 'Soto La Marina'
$ws.Range("B1033").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B1035").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B1037").Value = 'San Pablo Del Monte'
$ws.Range("B1039").Value = 'Tepetitla De Lardizábal'
$ws.Range("B1049").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B1053").Value = 'Amatlán De Los Reyes'
$ws.Range("B1060").Value = 'Boca Del Río'
$ws.Range("B1062").Value = 'Camarón De Tejeda'
$ws.Range("B1065").Value = 'Castillo De Teayo'
$ws.Range("B1067").Value = 'Cazones De Herrera'
$ws.Range("B1081").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1082").Value = 'Cosautlán De Carvajal'
$ws.Range("B1096").Value = 'Hueyapan De Ocampo'
$ws.Range("B1097").Value = 'Ignacio De La Llave'
$ws.Range("B1100").Value = 'Ixhuacán De Los Reyes'
$ws.Range("B1101").Value = 'Ixhuatlán De Madero'
$ws.Range("B1102").Value = 'Ixhuatlán Del Café'
$ws.Range("B1110").Value = 'Juchique De Ferrer'
$ws.Range("B1113").Value = 'Las Vigas De Ramírez'
$ws.Range("B1114").Value = 'Lerdo De Tejada'
$ws.Range("B1118").Value = 'Martínez De La Torre'
$ws.Range("B1120").Value = 'Medellín De Bravo'
$ws.Range("B1125").Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range("B1136").Value = 'Paso De Ovejas'
$ws.Range("B1137").Value = 'Paso Del Macho'
$ws.Range("B1140").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1149").Value = 'Sayula De Alemán'
$ws.Range("B1151").Value = 'Soledad De Doblado'
$ws.Range("B1155").Value = 'Tatahuicapan De Juárez'
$ws.Range("B1179").Value = 'Vega De Alatorre'
$ws.Range("B1189").Value = 'Zozocolco De Hidalgo'
$ws.Range("B1201").Value = 'Jiménez Del Teul'
$ws.Range("B1204").Value = 'Mezquital Del Oro'
$ws.Range("B1206").Value = 'Nochistlán De Mejía'
$ws.Range("B1207").Value = 'Noria De Ángeles'
$ws.Range("B1213").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1216").Value = 'Villa De Cos'

# Remove trailing metadata/footer rows (1224:1228); row 1222 dimension becomes A1:D1222
$ws.Rows("1224:1228").Delete()
